$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 102.4387563333333
$ws.Cells.Item(2, 8).Value = 307.316269
$ws.Cells.Item(2, 9).Value = 0.01835346024518914
$ws.Cells.Item(2, 10).Value = 0.0184635933110779
$ws.Cells.Item(2, 13).Value = 2.906846333333333
$ws.Cells.Item(2, 14).Value = 8.720538999999999
$ws.Cells.Item(2, 15).Value = 0.005520525738044089
$ws.Cells.Item(2, 16).Value = 0.005624540846623205
$ws.Cells.Item(2, 17).Value = 297.7737232387767
$ws.Cells.Item(2, 18).Value = 2679.96350914899
$ws.Cells.Item(2, 19).Value = 0.0001013207496657356
$ws.Cells.Item(2, 20).Value = 0.0001038492347535966
$ws.Cells.Item(3, 7).Value = 102.4387563333333
$ws.Cells.Item(3, 8).Value = 307.316269
$ws.Cells.Item(3, 9).Value = 0.01835346024518914
$ws.Cells.Item(3, 10).Value = 0.0184635933110779
$ws.Cells.Item(3, 15).Value = 0.3528665483720876
$ws.Cells.Item(3, 16).Value = 0.3595150912979765
$ws.Cells.Item(3, 17).Value = 19033.40205282696
$ws.Cells.Item(3, 18).Value = 171300.6184754427
$ws.Cells.Item(3, 19).Value = 0.00647632216740422
$ws.Cells.Item(3, 20).Value = 0.00663794043492088
$ws.Cells.Item(4, 7).Value = 102.4387563333333
$ws.Cells.Item(4, 8).Value = 307.316269
$ws.Cells.Item(4, 9).Value = 0.01835346024518914
$ws.Cells.Item(4, 10).Value = 0.0184635933110779
$ws.Cells.Item(4, 13).Value = 137.0717086666666
$ws.Cells.Item(4, 14).Value = 411.2151259999999
$ws.Cells.Item(4, 15).Value = 0.2603191943704447
$ws.Cells.Item(4, 16).Value = 0.2652240042658267
$ws.Cells.Item(4, 17).Value = 14041.45536429832
$ws.Cells.Item(4, 18).Value = 126373.0982786849
$ws.Cells.Item(4, 19).Value = 0.00477775798493762
$ws.Cells.Item(4, 20).Value = 0.004896988151099815
$ws.Cells.Item(5, 7).Value = 102.4387563333333
$ws.Cells.Item(5, 8).Value = 307.316269
$ws.Cells.Item(5, 9).Value = 0.01835346024518914
$ws.Cells.Item(5, 10).Value = 0.0184635933110779
$ws.Cells.Item(5, 13).Value = 29.2127365
$ws.Cells.Item(5, 14).Value = 58.425473
$ws.Cells.Item(5, 15).Value = 0.05547925319534149
$ws.Cells.Item(5, 16).Value = 0.03768304451958546
$ws.Cells.Item(5, 17).Value = 2992.516396153373
$ws.Cells.Item(5, 18).Value = 17955.09837692024
$ws.Cells.Item(5, 19).Value = 0.001018236267953482
$ws.Cells.Item(5, 20).Value = 0.0006957644087328687
$ws.Cells.Item(6, 7).Value = 102.4387563333333
$ws.Cells.Item(6, 8).Value = 307.316269
$ws.Cells.Item(6, 9).Value = 0.01835346024518914
$ws.Cells.Item(6, 10).Value = 0.0184635933110779
$ws.Cells.Item(6, 13).Value = 171.5584106666666
$ws.Cells.Item(6, 14).Value = 514.6752319999999
$ws.Cells.Item(6, 15).Value = 0.3258144783240821
$ws.Cells.Item(6, 16).Value = 0.331953319069988
$ws.Cells.Item(6, 17).Value = 17574.2302272166
$ws.Cells.Item(6, 18).Value = 158168.0720449494
$ws.Cells.Item(6, 19).Value = 0.005979823075228079
$ws.Cells.Item(6, 20).Value = 0.006129051081570739
$ws.Cells.Item(7, 9).Value = 0.03153528182220144
$ws.Cells.Item(7, 10).Value = 0.03172451465483067
$ws.Cells.Item(7, 13).Value = 2.906846333333333
$ws.Cells.Item(7, 14).Value = 8.720538999999999
$ws.Cells.Item(7, 15).Value = 0.005520525738044089
$ws.Cells.Item(7, 16).Value = 0.005624540846623205
$ws.Cells.Item(7, 17).Value = 511.6407563550567
$ws.Cells.Item(7, 18).Value = 4604.76680719551
$ws.Cells.Item(7, 19).Value = 0.0001740913349559369
$ws.Cells.Item(7, 20).Value = 0.0001784358285153916
$ws.Cells.Item(8, 9).Value = 0.03153528182220144
$ws.Cells.Item(8, 10).Value = 0.03172451465483067
$ws.Cells.Item(8, 15).Value = 0.3528665483720876
$ws.Cells.Item(8, 16).Value = 0.3595150912979765
$ws.Cells.Item(8, 19).Value = 0.01112774604854126
$ws.Cells.Item(8, 20).Value = 0.01140544178251544
$ws.Cells.Item(9, 9).Value = 0.03153528182220144
$ws.Cells.Item(9, 10).Value = 0.03172451465483067
$ws.Cells.Item(9, 13).Value = 137.0717086666666
$ws.Cells.Item(9, 14).Value = 411.2151259999999
$ws.Cells.Item(9, 15).Value = 0.2603191943704447
$ws.Cells.Item(9, 16).Value = 0.2652240042658267
$ws.Cells.Item(9, 17).Value = 24126.30894618784
$ws.Cells.Item(9, 18).Value = 217136.7805156906
$ws.Cells.Item(9, 19).Value = 0.008209239158200406
$ws.Cells.Item(9, 20).Value = 0.008414102810144092
$ws.Cells.Item(10, 9).Value = 0.03153528182220144
$ws.Cells.Item(10, 10).Value = 0.03172451465483067
$ws.Cells.Item(10, 13).Value = 29.2127365
$ws.Cells.Item(10, 14).Value = 58.425473
$ws.Cells.Item(10, 15).Value = 0.05547925319534149
$ws.Cells.Item(10, 16).Value = 0.03768304451958546
$ws.Cells.Item(10, 17).Value = 5141.801417800313
$ws.Cells.Item(10, 18).Value = 30850.80850680188
$ws.Cells.Item(10, 19).Value = 0.001749553884800363
$ws.Cells.Item(10, 20).Value = 0.001195476298100225
$ws.Cells.Item(11, 9).Value = 0.03153528182220144
$ws.Cells.Item(11, 10).Value = 0.03172451465483067
$ws.Cells.Item(11, 13).Value = 171.5584106666666
$ws.Cells.Item(11, 14).Value = 514.6752319999999
$ws.Cells.Item(11, 15).Value = 0.3258144783240821
$ws.Cells.Item(11, 16).Value = 0.331953319069988
$ws.Cells.Item(11, 17).Value = 30196.39324790524
$ws.Cells.Item(11, 18).Value = 271767.5392311472
$ws.Cells.Item(11, 19).Value = 0.01027465139570347
$ws.Cells.Item(11, 20).Value = 0.01053105793555552
$ws.Cells.Item(12, 7).Value = 2105.314697333333
$ws.Cells.Item(12, 8).Value = 6315.944092
$ws.Cells.Item(12, 9).Value = 0.3771991283785865
$ws.Cells.Item(12, 10).Value = 0.3794625760284536
$ws.Cells.Item(12, 13).Value = 2.906846333333333
$ws.Cells.Item(12, 14).Value = 8.720538999999999
$ws.Cells.Item(12, 15).Value = 0.005520525738044089
$ws.Cells.Item(12, 16).Value = 0.005624540846623205
$ws.Cells.Item(12, 17).Value = 6119.826308456175
$ws.Cells.Item(12, 18).Value = 55078.43677610558
$ws.Cells.Item(12, 19).Value = 0.002082337496581783
$ws.Cells.Item(12, 20).Value = 0.002134302758636901
$ws.Cells.Item(13, 7).Value = 2105.314697333333
$ws.Cells.Item(13, 8).Value = 6315.944092
$ws.Cells.Item(13, 9).Value = 0.3771991283785865
$ws.Cells.Item(13, 10).Value = 0.3794625760284536
$ws.Cells.Item(13, 15).Value = 0.3528665483720876
$ws.Cells.Item(13, 16).Value = 0.3595150912979765
$ws.Cells.Item(13, 17).Value = 391173.2484498344
$ws.Cells.Item(13, 18).Value = 3520559.236048509
$ws.Cells.Item(13, 19).Value = 0.1331009544799118
$ws.Cells.Item(13, 20).Value = 0.1364225226650349
$ws.Cells.Item(14, 7).Value = 2105.314697333333
$ws.Cells.Item(14, 8).Value = 6315.944092
$ws.Cells.Item(14, 9).Value = 0.3771991283785865
$ws.Cells.Item(14, 10).Value = 0.3794625760284536
$ws.Cells.Item(14, 13).Value = 137.0717086666666
$ws.Cells.Item(14, 14).Value = 411.2151259999999
$ws.Cells.Item(14, 15).Value = 0.2603191943704447
$ws.Cells.Item(14, 16).Value = 0.2652240042658267
$ws.Cells.Item(14, 17).Value = 288579.0828445261
$ws.Cells.Item(14, 18).Value = 2597211.745600735
$ws.Cells.Item(14, 19).Value = 0.09819217321674756
$ws.Cells.Item(14, 20).Value = 0.1006425838832922
$ws.Cells.Item(15, 7).Value = 2105.314697333333
$ws.Cells.Item(15, 8).Value = 6315.944092
$ws.Cells.Item(15, 9).Value = 0.3771991283785865
$ws.Cells.Item(15, 10).Value = 0.3794625760284536
$ws.Cells.Item(15, 13).Value = 29.2127365
$ws.Cells.Item(15, 14).Value = 58.425473
$ws.Cells.Item(15, 15).Value = 0.05547925319534149
$ws.Cells.Item(15, 16).Value = 0.03768304451958546
$ws.Cells.Item(15, 17).Value = 61502.00350277591
$ws.Cells.Item(15, 18).Value = 369012.0210166555
$ws.Cells.Item(15, 19).Value = 0.02092672594837772
$ws.Cells.Item(15, 20).Value = 0.0142993051459968
$ws.Cells.Item(16, 7).Value = 2105.314697333333
$ws.Cells.Item(16, 8).Value = 6315.944092
$ws.Cells.Item(16, 9).Value = 0.3771991283785865
$ws.Cells.Item(16, 10).Value = 0.3794625760284536
$ws.Cells.Item(16, 13).Value = 171.5584106666666
$ws.Cells.Item(16, 14).Value = 514.6752319999999
$ws.Cells.Item(16, 15).Value = 0.3258144783240821
$ws.Cells.Item(16, 16).Value = 0.331953319069988
$ws.Cells.Item(16, 17).Value = 361184.443427681
$ws.Cells.Item(16, 18).Value = 3250659.990849129
$ws.Cells.Item(16, 19).Value = 0.1228969372369676
$ws.Cells.Item(16, 20).Value = 0.1259638615754928
$ws.Cells.Item(17, 7).Value = 99.8778305
$ws.Cells.Item(17, 8).Value = 199.755661
$ws.Cells.Item(17, 9).Value = 0.01789463145660039
$ws.Cells.Item(17, 10).Value = 0.01200134082810157
$ws.Cells.Item(17, 13).Value = 2.906846333333333
$ws.Cells.Item(17, 14).Value = 8.720538999999999
$ws.Cells.Item(17, 15).Value = 0.005520525738044089
$ws.Cells.Item(17, 16).Value = 0.005624540846623205
$ws.Cells.Item(17, 17).Value = 290.3295053702131
$ws.Cells.Item(17, 18).Value = 1741.977032221279
$ws.Cells.Item(17, 19).Value = 0.00009878777352897585926471552
$ws.Cells.Item(17, 20).Value = 0.000067502031701904046594002462
$ws.Cells.Item(18, 7).Value = 99.8778305
$ws.Cells.Item(18, 8).Value = 199.755661
$ws.Cells.Item(18, 9).Value = 0.01789463145660039
$ws.Cells.Item(18, 10).Value = 0.01200134082810157
$ws.Cells.Item(18, 15).Value = 0.3528665483720876
$ws.Cells.Item(18, 16).Value = 0.3595150912979765
$ws.Cells.Item(18, 17).Value = 18557.5750049595
$ws.Cells.Item(18, 18).Value = 111345.450029757
$ws.Cells.Item(18, 19).Value = 0.006314416836481164
$ws.Cells.Item(18, 20).Value = 0.004314663143513069
$ws.Cells.Item(19, 7).Value = 99.8778305
$ws.Cells.Item(19, 8).Value = 199.755661
$ws.Cells.Item(19, 9).Value = 0.01789463145660039
$ws.Cells.Item(19, 10).Value = 0.01200134082810157
$ws.Cells.Item(19, 13).Value = 137.0717086666666
$ws.Cells.Item(19, 14).Value = 411.2151259999999
$ws.Cells.Item(19, 15).Value = 0.2603191943704447
$ws.Cells.Item(19, 16).Value = 0.2652240042658267
$ws.Cells.Item(19, 17).Value = 13690.42488455471
$ws.Cells.Item(19, 18).Value = 82142.54930732827
$ws.Cells.Item(19, 19).Value = 0.004658316044338231
$ws.Cells.Item(19, 20).Value = 0.003183043670988051
$ws.Cells.Item(20, 7).Value = 99.8778305
$ws.Cells.Item(20, 8).Value = 199.755661
$ws.Cells.Item(20, 9).Value = 0.01789463145660039
$ws.Cells.Item(20, 10).Value = 0.01200134082810157
$ws.Cells.Item(20, 13).Value = 29.2127365
$ws.Cells.Item(20, 14).Value = 58.425473
$ws.Cells.Item(20, 15).Value = 0.05547925319534149
$ws.Cells.Item(20, 16).Value = 0.03768304451958546
$ws.Cells.Item(20, 17).Value = 2917.704744588163
$ws.Cells.Item(20, 18).Value = 11670.81897835265
$ws.Cells.Item(20, 19).Value = 0.0009927807894180556
$ws.Cells.Item(20, 20).Value = 0.0004522470607200701
$ws.Cells.Item(21, 7).Value = 99.8778305
$ws.Cells.Item(21, 8).Value = 199.755661
$ws.Cells.Item(21, 9).Value = 0.01789463145660039
$ws.Cells.Item(21, 10).Value = 0.01200134082810157
$ws.Cells.Item(21, 13).Value = 171.5584106666666
$ws.Cells.Item(21, 14).Value = 514.6752319999999
$ws.Cells.Item(21, 15).Value = 0.3258144783240821
$ws.Cells.Item(21, 16).Value = 0.331953319069988
$ws.Cells.Item(21, 17).Value = 17134.88186141472
$ws.Cells.Item(21, 18).Value = 102809.2911684883
$ws.Cells.Item(21, 19).Value = 0.005830330012833967
$ws.Cells.Item(21, 20).Value = 0.003983884921178476
$ws.Cells.Item(22, 7).Value = 3097.797444666667
$ws.Cells.Item(22, 8).Value = 9293.392334
$ws.Cells.Item(22, 9).Value = 0.5550174980974226
$ws.Cells.Item(22, 10).Value = 0.5583479751775362
$ws.Cells.Item(22, 13).Value = 2.906846333333333
$ws.Cells.Item(22, 14).Value = 8.720538999999999
$ws.Cells.Item(22, 15).Value = 0.005520525738044089
$ws.Cells.Item(22, 16).Value = 0.005624540846623205
$ws.Cells.Item(22, 17).Value = 9004.821143438669
$ws.Cells.Item(22, 18).Value = 81043.39029094801
$ws.Cells.Item(22, 19).Value = 0.003063988383311658
$ws.Cells.Item(22, 20).Value = 0.003140450993015412
$ws.Cells.Item(23, 7).Value = 3097.797444666667
$ws.Cells.Item(23, 8).Value = 9293.392334
$ws.Cells.Item(23, 9).Value = 0.5550174980974226
$ws.Cells.Item(23, 10).Value = 0.5583479751775362
$ws.Cells.Item(23, 15).Value = 0.3528665483720876
$ws.Cells.Item(23, 16).Value = 0.3595150912979765
$ws.Cells.Item(23, 17).Value = 575579.2666078539
$ws.Cells.Item(23, 18).Value = 5180213.399470686
$ws.Cells.Item(23, 19).Value = 0.1958471088397492
$ws.Cells.Item(23, 20).Value = 0.2007345232719923
$ws.Cells.Item(24, 7).Value = 3097.797444666667
$ws.Cells.Item(24, 8).Value = 9293.392334
$ws.Cells.Item(24, 9).Value = 0.5550174980974226
$ws.Cells.Item(24, 10).Value = 0.5583479751775362
$ws.Cells.Item(24, 13).Value = 137.0717086666666
$ws.Cells.Item(24, 14).Value = 411.2151259999999
$ws.Cells.Item(24, 15).Value = 0.2603191943704447
$ws.Cells.Item(24, 16).Value = 0.2652240042658267
$ws.Cells.Item(24, 17).Value = 424620.3888436937
$ws.Cells.Item(24, 18).Value = 3821583.499593243
$ws.Cells.Item(24, 19).Value = 0.1444817079662208
$ws.Cells.Item(24, 20).Value = 0.1480872857503026
$ws.Cells.Item(25, 7).Value = 3097.797444666667
$ws.Cells.Item(25, 8).Value = 9293.392334
$ws.Cells.Item(25, 9).Value = 0.5550174980974226
$ws.Cells.Item(25, 10).Value = 0.5583479751775362
$ws.Cells.Item(25, 13).Value = 29.2127365
$ws.Cells.Item(25, 14).Value = 58.425473
$ws.Cells.Item(25, 15).Value = 0.05547925319534149
$ws.Cells.Item(25, 16).Value = 0.03768304451958546
$ws.Cells.Item(25, 17).Value = 90495.14048142066
$ws.Cells.Item(25, 18).Value = 542970.8428885239
$ws.Cells.Item(25, 19).Value = 0.03079195630479187
$ws.Cells.Item(25, 20).Value = 0.02104025160603549
$ws.Cells.Item(26, 7).Value = 3097.797444666667
$ws.Cells.Item(26, 8).Value = 9293.392334
$ws.Cells.Item(26, 9).Value = 0.5550174980974226
$ws.Cells.Item(26, 10).Value = 0.5583479751775362
$ws.Cells.Item(26, 13).Value = 171.5584106666666
$ws.Cells.Item(26, 14).Value = 514.6752319999999
$ws.Cells.Item(26, 15).Value = 0.3258144783240821
$ws.Cells.Item(26, 16).Value = 0.331953319069988
$ws.Cells.Item(26, 17).Value = 531453.2061742746
$ws.Cells.Item(26, 18).Value = 4783078.855568471
$ws.Cells.Item(26, 19).Value = 0.180832736603349
$ws.Cells.Item(26, 20).Value = 0.18534546355619042
